$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------- Row 2 ----------
Set-TextValue $ws.Range("A2") "281474991265672-1748989088169"
$ws.Range("B2").Value = "Harsh Brake"
$ws.Range("C2").Value = "2025-06-03T16:18:08.169"
Set-TextValue $ws.Range("D2") "281474991265672"
Set-TextValue $ws.Range("E2") "116"
Set-TextValue $ws.Range("F2") "52215735"
$ws.Range("G2").Value = "KEVIN DE LA O"
$ws.Range("H2").Value = 20.645476299
$ws.Range("I2").Value = -103.34981778
$ws.Range("J2").Value = 0.7115712761878967
$ws.Range("K2").Value = "No video URL"
$ws.Range("L2").Value = "No video URL"

# ---------- Row 3 ----------
Set-TextValue $ws.Range("A3") "281474990867465-1748977464722"
$ws.Range("B3").Value = "No Seat Belt"
$ws.Range("C3").Value = "2025-06-03T13:04:24.722"
Set-TextValue $ws.Range("D3") "281474990867465"
Set-TextValue $ws.Range("E3") "124"
Set-TextValue $ws.Range("F3") "52215846"
$ws.Range("G3").Value = "ALVARO ZAPATA"
$ws.Range("H3").Value = 20.74507293
$ws.Range("I3").Value = -105.4243703
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1748977462222/8CVckh4ELA-camera-video-segment-driver-1748977464722.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCV4U56KD%2F20250604%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250604T150921Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEFcaCXVzLXdlc3QtMiJHMEUCIQCBcdkIHf0VU8WJbQqltdniPAWqv%2FAzaTLnfiU6VFnjwgIgHJR75YNhoVHgvBTiVX3SqtJr2lcnmjg2Lj24GxSjy%2F0q3QMIMBAEGgw3ODEyMDQ5NDIyNDQiDM48xlDf6Y6LTW2ZpCq6AzjkF%2BAvhWDuhXJQQcSuZkye%2FQNAZHhdM1%2B1nMub5iill5nH%2BvYD%2FiarxyHuIPcd0eFJ9iU9kNYfQOexBaMLqnw2Zu4%2BSta6ZY6cgbGWxYc6GeicyvVHnr20FR%2FMSNFLy58mi%2BGhticoieFBCmGW0JOpnYH2TH7i36hLLCP0DFw5W4ufiqsKT7T39f2NCgHFifFLLGGfr3rlWynEvXFC4DkhQU64RfOfHt8AC6%2F9DtMRtjYP0238bPWF0TTcxZDrRQZo6%2FEcHFXQeJEngJKyNbm%2BrasboL8MDRKYigcNfNv34eXSWnUpjaGBPW%2B%2FFA2cnpOfzv1UlgZm%2FN88N6M8%2FDwJHWWjQBReEArMiNR5S%2F8MUkHbmlbdZi620wb7OxVuD%2FnacXhOvjB8KOE7yfqrU8u2m655CX2B01z2bJSs9c1OHOQwqFiDBmGgpp16N1QjLvo8PexvlYrhx5a3vUbWI6tSBKLIEFQpwL2vooLEYqWMaxIRBFls5o%2F589O9DnlmLLMwnfp6rvaFepoqvaNiARucVE9l%2FBg38S%2BzGTmJrOLdpNPX9sADwKQ6zIaNEgnwRgWt7jn2K63OOOgw4rKBwgY6pQEaSxRXmI6l%2BUSFZWy44GdVJvCOX%2FxXHFpFDDvTsgv%2BFWVDf7%2BbgGNJRpC6%2Bu%2FgposgVJj3EeTmqn6wkaTHi4glf23C2Na%2FguYmZ8fcw9pA%2Bd9PhzLWyA8yLUMijcGQo7IXC%2FruGJgTgo7bxGYN%2B0uBYg08xI4N0FHaIt747s87SmoaJKsoDUnntNAwmh6eEiYwjuxu9nlppi5a%2Bc7HV01tyxHqtgk%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2023%3A09%3A21%20GMT&X-Amz-Signature=9da0fdbaf0686a810990eccc1772de78e481e33ab7852fc8e4c3414892f04087"
$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1748977462222/f11BxoPt2w-camera-video-segment-1748977464722.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCV4U56KD%2F20250604%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250604T150921Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEFcaCXVzLXdlc3QtMiJHMEUCIQCBcdkIHf0VU8WJbQqltdniPAWqv%2FAzaTLnfiU6VFnjwgIgHJR75YNhoVHgvBTiVX3SqtJr2lcnmjg2Lj24GxSjy%2F0q3QMIMBAEGgw3ODEyMDQ5NDIyNDQiDM48xlDf6Y6LTW2ZpCq6AzjkF%2BAvhWDuhXJQQcSuZkye%2FQNAZHhdM1%2B1nMub5iill5nH%2BvYD%2FiarxyHuIPcd0eFJ9iU9kNYfQOexBaMLqnw2Zu4%2BSta6ZY6cgbGWxYc6GeicyvVHnr20FR%2FMSNFLy58mi%2BGhticoieFBCmGW0JOpnYH2TH7i36hLLCP0DFw5W4ufiqsKT7T39f2NCgHFifFLLGGfr3rlWynEvXFC4DkhQU64RfOfHt8AC6%2F9DtMRtjYP0238bPWF0TTcxZDrRQZo6%2FEcHFXQeJEngJKyNbm%2BrasboL8MDRKYigcNfNv34eXSWnUpjaGBPW%2B%2FFA2cnpOfzv1UlgZm%2FN88N6M8%2FDwJHWWjQBReEArMiNR5S%2F8MUkHbmlbdZi620wb7OxVuD%2FnacXhOvjB8KOE7yfqrU8u2m655CX2B01z2bJSs9c1OHOQwqFiDBmGgpp16N1QjLvo8PexvlYrhx5a3vUbWI6tSBKLIEFQpwL2vooLEYqWMaxIRBFls5o%2F589O9DnlmLLMwnfp6rvaFepoqvaNiARucVE9l%2FBg38S%2BzGTmJrOLdpNPX9sADwKQ6zIaNEgnwRgWt7jn2K63OOOgw4rKBwgY6pQEaSxRXmI6l%2BUSFZWy44GdVJvCOX%2FxXHFpFDDvTsgv%2BFWVDf7%2BbgGNJRpC6%2Bu%2FgposgVJj3EeTmqn6wkaTHi4glf23C2Na%2FguYmZ8fcw9pA%2Bd9PhzLWyA8yLUMijcGQo7IXC%2FruGJgTgo7bxGYN%2B0uBYg08xI4N0FHaIt747s87SmoaJKsoDUnntNAwmh6eEiYwjuxu9nlppi5a%2Bc7HV01tyxHqtgk%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2004%20Jun%202025%2023%3A09%3A21%20GMT&X-Amz-Signature=56c16d32bff60e41cb04f8a7fb792365957790dd55051390ded9c7db93dca867"

# ---------- Row 4 (new row) ----------
Set-TextValue $ws.Range("A4") "281474991395157-1748969368760"
$ws.Range("B4").Value = "Harsh Brake"
$ws.Range("C4").Value = "2025-06-03T10:49:28.760"
Set-TextValue $ws.Range("D4") "281474991395157"
Set-TextValue $ws.Range("E4") "126"
$ws.Range("F4").Value = "No driver ID"
$ws.Range("G4").Value = "No driver name"
$ws.Range("H4").Value = 20.62094963
$ws.Range("I4").Value = -103.422907179
$ws.Range("J4").Value = 0.7334661483764648
$ws.Range("K4").Value = "No video URL"
$ws.Range("L4").Value = "No video URL"
